# Update quality check report
# Insert two new dictionary rows (crfs-t09a2-i2_1b / crfs-t09a2-i2_1o) just
# above the existing row that used to be row 163 ("crfs-t09a2-j2_1"),
# pushing all the subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 163 (existing rows 163:205 shift down to 165:207).
# Excel's default Insert copies formatting from the row above, which gives us
# the same style ids (17/14/2/2/20/22/3/17) used by the surrounding rows.
$ws.Rows.Item(163).Insert()
$ws.Rows.Item(163).Insert()

# New row 163: crfs-t09a2-i2_1b
$ws.Cells.Item(163, 1).Value = "crfs-t09a2-i2_1b"
$ws.Cells.Item(163, 2).Value = "ref_facility"
$ws.Cells.Item(163, 3).Value = 0
$ws.Cells.Item(163, 4).Value = 0
$ws.Cells.Item(163, 5).Value = 1
$ws.Cells.Item(163, 6).Value = 1
$ws.Cells.Item(163, 7).Value = 1
$ws.Cells.Item(163, 8).Value = "i2_1b"

# New row 164: crfs-t09a2-i2_1o
$ws.Cells.Item(164, 1).Value = "crfs-t09a2-i2_1o"
$ws.Cells.Item(164, 2).Value = "ref_facility_oth"
$ws.Cells.Item(164, 3).Value = 0
$ws.Cells.Item(164, 4).Value = 0
$ws.Cells.Item(164, 5).Value = 1
$ws.Cells.Item(164, 6).Value = 1
$ws.Cells.Item(164, 7).Value = 1
$ws.Cells.Item(164, 8).Value = "i2_1o"

# Reflect the editor's final scroll position / active selection.
$excel.ActiveWindow.ScrollRow = 136
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J153").Select()
